$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# productname changed to the "-1st" suffixed test case name (dedup to remove
# test-case inter-dependency) on both the input and output sheets.
$newProductName = "2627-MS-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-AMT-VAR-INST-FIX-INST-AMT-MORE-1st"
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# shortname changed from the numeric 2627 to the text value "262d"
$wsInput.Range("B2").Value = "262d"

# Move the active selection from B15 to B3
$wsInput.Range("B3").Select() | Out-Null
